$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.309.17'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '1.620.43'
$ws.Range('E3').Value = '  +1.81%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('E7').Value = '  +0.83%  '
$ws.Range('E8').Value = '  +0.34%  '
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.75'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0815'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('D12').Value = '1.846.31'
$ws.Range('E12').Value = '  +1.81%  '
$ws.Range('D13').Value = '1.615.64'
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('E14').Value = '  +0.45%  '
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D16').Value = '26.313.38'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.29'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.50%  '
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '201.63'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('E21').Value = '  +1.12%  '
$ws.Range('E22').Value = '  +1.43%  '
$ws.Range('E23').Value = '  +0.76%  '
$ws.Range('E24').Value = '  -4.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.72'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.07%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('E29').Value = '  +1.43%  '
$ws.Range('E30').Value = '  +8.74%  '
$ws.Range('E31').Value = '  +0.61%  '
$ws.Range('E32').Value = '  +1.87%  '
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.49'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.49%  '
$ws.Range('E35').Value = '  +2.61%  '
$ws.Range('D36').Value = '1.178.97'
$ws.Range('E36').Value = '  +4.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0163'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('E38').Value = '  +2.62%  '
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('E41').Value = '  +1.28%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.35'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.28%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.785'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.60%  '
$ws.Range('D44').Value = '1.757.78'
$ws.Range('E44').Value = '  +1.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.62'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.53'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '53.77'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('E50').Value = '  -0.22%  '
